$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1169.2632
$ws.Range("I19").Value = 979.2
$ws.Range("J19").Value = 1380.4445
$ws.Range("K19").Value = 979.2
$ws.Range("L19").Value = 1380.4445
$ws.Range("M19").Value = -804.2
$ws.Range("N19").Value = -1730.4445

$ws.Range("H116").Value = 2021
$ws.Range("I116").Value = 2475
$ws.Range("J116").Value = 1761.5714
$ws.Range("K116").Value = 2475
$ws.Range("L116").Value = 1761.5714
$ws.Range("M116").Value = 967
$ws.Range("N116").Value = -8645.571400000001

$ws.Range("H131").Value = 4201.4907
$ws.Range("I131").Value = 1028.9
$ws.Range("J131").Value = 4939.3022
$ws.Range("K131").Value = 3086.7
$ws.Range("L131").Value = 14817.9066
$ws.Range("M131").Value = 1953.3
$ws.Range("N131").Value = -24897.9066

$ws.Range("H138").Value = 3684.7734
$ws.Range("I138").Value = 4513
$ws.Range("J138").Value = 3599.5146
$ws.Range("K138").Value = 13539
$ws.Range("L138").Value = 10798.5438
$ws.Range("M138").Value = -8399
$ws.Range("N138").Value = -21078.5438

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 42501.5
$ws.Range("I6").Value = 50002
$ws.Range("J6").Value = 20000
$ws.Range("K6").Value = 50002
$ws.Range("L6").Value = 20000
$ws.Range("M6").Value = -49829
$ws.Range("N6").Value = -20346

$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H23").Value = 34668.668
$ws.Range("I23").Value = 80006
$ws.Range("J23").Value = 12000
$ws.Range("K23").Value = 80006
$ws.Range("L23").Value = 12000
$ws.Range("M23").Value = -79747
$ws.Range("N23").Value = -12518

$ws.Range("H37").Value = 10900
$ws.Range("I37").Value = 6800
$ws.Range("J37").Value = 15000
$ws.Range("K37").Value = 6800
$ws.Range("L37").Value = 15000
$ws.Range("M37").Value = -6527
$ws.Range("N37").Value = -15546

$ws.Range("H63").Value = 3100
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 3100
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 3100
$ws.Range("N63").Value = -4472

$ws.Range("H66").Value = 3100
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 3100
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 15500
$ws.Range("N66").Value = -22364

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H35").Value = 19558
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 19558
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 19558
$ws.Range("N35").Value = -20178

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H80").Value = 1099.3939
$ws.Range("I80").Value = 640.0714
$ws.Range("J80").Value = 1437.8422
$ws.Range("K80").Value = 640.0714
$ws.Range("L80").Value = 1437.8422
$ws.Range("M80").Value = 357.9286
$ws.Range("N80").Value = -3433.8422

$ws.Range("H82").Value = 15810.4
$ws.Range("I82").Value = 3164
$ws.Range("J82").Value = 34780
$ws.Range("K82").Value = 3164
$ws.Range("L82").Value = 34780
$ws.Range("M82").Value = -2781
$ws.Range("N82").Value = -35546

$ws.Range("H83").Value = 1099.3939
$ws.Range("I83").Value = 640.0714
$ws.Range("J83").Value = 1437.8422
$ws.Range("K83").Value = 3200.357
$ws.Range("L83").Value = 7189.211
$ws.Range("M83").Value = 1791.643
$ws.Range("N83").Value = -17173.211

$ws.Range("H85").Value = 15810.4
$ws.Range("I85").Value = 3164
$ws.Range("J85").Value = 34780
$ws.Range("K85").Value = 3164
$ws.Range("L85").Value = 34780
$ws.Range("M85").Value = -1838
$ws.Range("N85").Value = -37432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1390.9
$ws.Range("I105").Value = 1391.125
$ws.Range("J105").Value = 1390
$ws.Range("K105").Value = 1391.125
$ws.Range("L105").Value = 1390
$ws.Range("M105").Value = 355.875
$ws.Range("N105").Value = -4884

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 73678.57000000001
$ws.Range("I70").Value = 167918.33
$ws.Range("J70").Value = 2998.75
$ws.Range("K70").Value = 503754.99
$ws.Range("L70").Value = 8996.25
$ws.Range("M70").Value = -503439.99
$ws.Range("N70").Value = -9626.25

$ws.Range("H73").Value = 73678.57000000001
$ws.Range("I73").Value = 167918.33
$ws.Range("J73").Value = 2998.75
$ws.Range("K73").Value = 503754.99
$ws.Range("L73").Value = 8996.25
$ws.Range("M73").Value = -502662.99
$ws.Range("N73").Value = -11180.25

$ws.Range("H131").Value = 713752.4399999999
$ws.Range("I131").Value = 586
$ws.Range("J131").Value = 830664.9399999999
$ws.Range("K131").Value = 1758
$ws.Range("L131").Value = 2491994.82
$ws.Range("M131").Value = 3282
$ws.Range("N131").Value = -2502074.82

$ws.Range("H132").Value = 2079.4736
$ws.Range("I132").Value = 1040
$ws.Range("J132").Value = 2450.7144
$ws.Range("K132").Value = 9360
$ws.Range("L132").Value = 22056.4296
$ws.Range("M132").Value = -6830
$ws.Range("N132").Value = -27116.4296

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1221.5333
$ws.Range("I113").Value = 834.55554
$ws.Range("J113").Value = 1802
$ws.Range("K113").Value = 834.55554
$ws.Range("L113").Value = 1802
$ws.Range("M113").Value = 1335.44446
$ws.Range("N113").Value = -6142

$ws.Range("H132").Value = 2827.3022
$ws.Range("I132").Value = 2181.2222
$ws.Range("J132").Value = 3917.5625
$ws.Range("K132").Value = 6543.6666
$ws.Range("L132").Value = 11752.6875
$ws.Range("M132").Value = -4013.6666
$ws.Range("N132").Value = -16812.6875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 37000
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 37000
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 37000
$ws.Range("N36").Value = -38124

$ws.Range("H128").Value = 48996
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 48996
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 48996
$ws.Range("N128").Value = -58956

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 38555.332
$ws.Range("I26").Value = 6666
$ws.Range("J26").Value = 54500
$ws.Range("K26").Value = 6666
$ws.Range("L26").Value = 54500
$ws.Range("M26").Value = -6373
$ws.Range("N26").Value = -55086

$ws.Range("H128").Value = 51838.332
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 51838.332
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 51838.332
$ws.Range("N128").Value = -61798.332
